# Updated the python code to handle the file does not exist error.
#
# The data that used to be generated by the "Literacy" data-source query
# is no longer available (the source file couldn't be found), so the
# connection-level metadata on the "Data Source" sheet is wiped back to an
# empty placeholder table, and the captured LOAD script on the
# "Data Source Query" sheet is replaced with "NA" markers.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Data Source" sheet / table2: collapse the 6-column x 1-row table
# (Dashboard Name / Connection Name / Connection ID / Data Source /
# Data Source Type / Table Name) down to a single placeholder column
# named "None", with no data rows left at all.
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data Source")
$loData = $wsData.ListObjects.Item(1)

# Shrink the table to a single header cell and rename that column to "None"
# (renaming a ListColumn is driven by its header-cell text).
$loData.Resize($wsData.Range("A1:A1"))
$wsData.Range("A1").Value = "None"

# Nothing is actually left behind in the sheet data -- clear the header
# placeholder text and the old data row, then physically drop the now
# unused columns B:F plus column A itself so no stray cell values remain.
$wsData.Range("A1").ClearContents()
$wsData.Range("A2").ClearContents()
$wsData.Range("A:F").Delete()

# ------------------------------------------------------------------
# "Data Source Query" sheet: the captured LOAD script for
# [LiteracyRate] is gone -- replace both the table-name and the
# source-query cells with "NA", and shrink the columns that used to
# hold that long query text.
# ------------------------------------------------------------------
$wsQuery = $wb.Worksheets.Item("Data Source Query")
$wsQuery.Columns.Item(2).ColumnWidth = 13.5
$wsQuery.Columns.Item(3).ColumnWidth = 16
$wsQuery.Range("B2").Value = "NA"
$wsQuery.Range("C2").Value = "NA"
